# Update the Implementation Guide metadata / generated spreadsheet to match
# the new canonical URL scheme (subdomain instead of path segment) and the
# refreshed generation timestamp.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("B2").Value = "https://molic-avc.gabriellesantosleandro.com/StructureDefinition/molicavc-ethnicity-extension"
$wsMeta.Range("B8").Value = "2023-08-16T00:27:03-03:00"

# --- Elements sheet -----------------------------------------------------
$wsElements = $wb.Worksheets.Item("Elements")

$wsElements.Range("R4").Value = "https://molic-avc.gabriellesantosleandro.com/StructureDefinition/molicavc-ethnicity-extension"
$wsElements.Range("Z5").Value = "https://molic-avc.gabriellesantosleandro.com/ValueSet/molicavc-ethnicity-valueset"

# The generated "best fit" width of column Z shrinks now that the URL text
# is shorter, matching the authoring tool's automatic column sizing
# (target best-fit width ~75.97 characters; 75.2 resolves to the nearest
# width the host's column-width grid can represent).
$wsElements.Columns.Item(26).ColumnWidth = 75.2
